# Weekly refresh of Fruta / hortaliza price data.
# Updates "Fecha" (D) and price columns (J,K,L,M,P) for rows 2-22
# to reflect the latest weekly reshuffled dataset.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44229
$ws.Range("J2").Value = 1500
$ws.Range("K2").Value = 1400
$ws.Range("L2").Value = 1500
$ws.Range("M2").Value = 1450
$ws.Range("P2").Value = 1450
$ws.Range("D3").Value = 44673
$ws.Range("J3").Value = 900
$ws.Range("K3").Value = 1300
$ws.Range("L3").Value = 1400
$ws.Range("M3").Value = 1350
$ws.Range("P3").Value = 1350
$ws.Range("D4").Value = 44687
$ws.Range("J4").Value = 1000
$ws.Range("K4").Value = 1200
$ws.Range("L4").Value = 1300
$ws.Range("M4").Value = 1250
$ws.Range("P4").Value = 1250
$ws.Range("D5").Value = 44175
$ws.Range("J5").Value = 1600
$ws.Range("K5").Value = 1000
$ws.Range("L5").Value = 1200
$ws.Range("M5").Value = 1100
$ws.Range("P5").Value = 1100
$ws.Range("D6").Value = 44550
$ws.Range("J6").Value = 1300
$ws.Range("K6").Value = 1000
$ws.Range("L6").Value = 1200
$ws.Range("M6").Value = 1100
$ws.Range("P6").Value = 1100
$ws.Range("D7").Value = 44607
$ws.Range("J7").Value = 900
$ws.Range("K7").Value = 1300
$ws.Range("L7").Value = 1400
$ws.Range("M7").Value = 1350
$ws.Range("P7").Value = 1350
$ws.Range("D8").Value = 44453
$ws.Range("J8").Value = 1000
$ws.Range("L8").Value = 900
$ws.Range("M8").Value = 850
$ws.Range("P8").Value = 850
$ws.Range("D9").Value = 44284
$ws.Range("J9").Value = 1500
$ws.Range("K9").Value = 800
$ws.Range("L9").Value = 850
$ws.Range("M9").Value = 825
$ws.Range("P9").Value = 825
$ws.Range("D10").Value = 44649
$ws.Range("J10").Value = 600
$ws.Range("D11").Value = 44407
$ws.Range("J11").Value = 1000
$ws.Range("K11").Value = 1200
$ws.Range("L11").Value = 1300
$ws.Range("M11").Value = 1250
$ws.Range("P11").Value = 1250
$ws.Range("D12").Value = 44784
$ws.Range("D13").Value = 44341
$ws.Range("J13").Value = 1300
$ws.Range("L13").Value = 1000
$ws.Range("M13").Value = 950
$ws.Range("P13").Value = 950
$ws.Range("D14").Value = 44243
$ws.Range("J14").Value = 1200
$ws.Range("D15").Value = 44656
$ws.Range("J15").Value = 1000
$ws.Range("K15").Value = 900
$ws.Range("L15").Value = 1000
$ws.Range("M15").Value = 950
$ws.Range("P15").Value = 950
$ws.Range("D16").Value = 44638
$ws.Range("J16").Value = 1000
$ws.Range("L16").Value = 950
$ws.Range("M16").Value = 925
$ws.Range("P16").Value = 925
$ws.Range("D17").Value = 44476
$ws.Range("K17").Value = 700
$ws.Range("L17").Value = 800
$ws.Range("M17").Value = 750
$ws.Range("P17").Value = 750
$ws.Range("D18").Value = 44291
$ws.Range("J18").Value = 1000
$ws.Range("K18").Value = 1000
$ws.Range("L18").Value = 1200
$ws.Range("M18").Value = 1100
$ws.Range("P18").Value = 1100
$ws.Range("D19").Value = 44449
$ws.Range("J19").Value = 1300
$ws.Range("K19").Value = 900
$ws.Range("L19").Value = 950
$ws.Range("M19").Value = 925
$ws.Range("P19").Value = 925
$ws.Range("D20").Value = 44442
$ws.Range("J20").Value = 1250
$ws.Range("K20").Value = 850
$ws.Range("L20").Value = 900
$ws.Range("M20").Value = 875
$ws.Range("P20").Value = 875
$ws.Range("D21").Value = 44455
$ws.Range("J21").Value = 1100
$ws.Range("K21").Value = 900
$ws.Range("L21").Value = 1000
$ws.Range("M21").Value = 950
$ws.Range("P21").Value = 950
$ws.Range("D22").Value = 44484
$ws.Range("J22").Value = 900
$ws.Range("K22").Value = 750
$ws.Range("L22").Value = 800
$ws.Range("M22").Value = 775
$ws.Range("P22").Value = 775